$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("G13:G18")
$range.Select()
$range.ClearContents()
